$d = $word.ActiveDocument

# --- Table cell corrections: drop trailing ":0" ---
$d.Content.Find.Execute("44:17:29:0", $true, $false, $false, $false, $false, $true, 1, $false, "44:17:29", 2)
$d.Content.Find.Execute("45:01:20:0", $true, $false, $false, $false, $false, $true, 1, $false, "45:01:20", 2)
$d.Content.Find.Execute("51:02:09:0", $true, $false, $false, $false, $false, $true, 1, $false, "51:02:09", 2)
$d.Content.Find.Execute("0:0:0:0", $true, $false, $false, $false, $false, $true, 1, $false, "0:0:0", 2)

# --- Re-type (self-replace) text spans so Word merges runs & drops proofErr markers ---
$d.Content.Find.Execute("/find/silver_gold/?span=7", $true, $false, $false, $false, $false, $true, 1, $false, "/find/silver_gold/?span=7", 2)

$d.Content.Find.Execute("/find", $true, $false, $false, $false, $false, $true, 1, $false, "/find", 2)

$d.Content.Find.Execute("NCODE everything but space [spaces are squenched] characters are ToLower(),variables . Also, ellipsis are always delimited with spaces.", $true, $false, $false, $false, $false, $true, 1, $false, "NCODE everything but space [spaces are squenched] characters are ToLower(),variables . Also, ellipsis are always delimited with spaces.", 2)

$d.Content.Find.Execute("All variables are sorted and ToLower() on the back-end to turn the search into a hashable normalized ID.", $true, $false, $false, $false, $false, $true, 1, $false, "All variables are sorted and ToLower() on the back-end to turn the search into a hashable normalized ID.", 2)

$d.Content.Find.Execute("/reference/Acts/17/find-quoted/silver_gold/?span=7&lexicon.search=dual", $true, $false, $false, $false, $false, $true, 1, $false, "/reference/Acts/17/find-quoted/silver_gold/?span=7&lexicon.search=dual", 2)

$d.Content.Find.Execute("/reference/44/17/?lexicon.display=av", $true, $false, $false, $false, $false, $true, 1, $false, "/reference/44/17/?lexicon.display=av", 2)

$d.Content.Find.Execute("We need also an API specific variable to indicate difference-highlighting [true/false ; 1/0]", $true, $false, $false, $false, $false, $true, 1, $false, "We need also an API specific variable to indicate difference-highlighting [true/false ; 1/0]", 2)
